$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -12
$ws.Range("F3").Value = -11
$ws.Range("F7").Value = 8
$ws.Range("E8").Value = -2
$ws.Range("F8").Value = -6
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 8
$ws.Range("F9").Value = -4
$ws.Range("F11").Value = -1
$ws.Range("F13").Value = -6
$ws.Range("F15").Value = -5
$ws.Range("F18").Value = -1
$ws.Range("F21").Value = -4
$ws.Range("F24").Value = -4
$ws.Range("F27").Value = -3
$ws.Range("F30").Value = -1
$ws.Range("F33").Value = 2
$ws.Range("F36").Value = -2
$ws.Range("F38").Value = 1
$ws.Range("F39").Value = -3
$ws.Range("F40").Value = 0
$ws.Range("F43").Value = -4
$ws.Range("F44").Value = -1
$ws.Range("F45").Value = -1
$ws.Range("F46").Value = -4
$ws.Range("F47").Value = 7
